$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.104.09'
$ws.Range("E2").Value = '  +4.47%  '
$ws.Range("D3").Value = '1.909.22'
$ws.Range("E3").Value = '  +5.32%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''251.27'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '''0.5090'
$ws.Range("E7").Value = '  +2.70%  '
$ws.Range("D8").Value = '''44.72'
$ws.Range("E8").Value = '  +3.18%  '
$ws.Range("D9").Value = '''0.2952'
$ws.Range("E9").Value = '  +5.84%  '
$ws.Range("D10").Value = '''0.06761'
$ws.Range("E10").Value = '  +5.56%  '
$ws.Range("D11").Value = '1.911.69'
$ws.Range("E11").Value = '  +5.47%  '
$ws.Range("D12").Value = '''17.23'
$ws.Range("E12").Value = '  +2.66%  '
$ws.Range("D13").Value = '''0.07361'
$ws.Range("E13").Value = '  +3.10%  '
$ws.Range("D14").Value = '''0.6884'
$ws.Range("E14").Value = '  +5.78%  '
$ws.Range("D15").Value = '''86.38'
$ws.Range("E15").Value = '  +2.95%  '
$ws.Range("D16").Value = '''4.868'
$ws.Range("E16").Value = '  +3.53%  '
$ws.Range("D17").Value = '30.115.52'
$ws.Range("E17").Value = '  +4.56%  '
$ws.Range("E18").Value = '  +9.23%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '''12.96'
$ws.Range("E20").Value = '  +5.83%  '
$ws.Range("D21").Value = '2.160.01'
$ws.Range("E21").Value = '  +5.25%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '''4.824'
$ws.Range("E23").Value = '  +4.63%  '
$ws.Range("D24").Value = '''5.734'
$ws.Range("E24").Value = '  +7.07%  '
$ws.Range("D25").Value = '''9.138'
$ws.Range("E25").Value = '  +2.65%  '
$ws.Range("D26").Value = '''146.86'
$ws.Range("E26").Value = '  +2.56%  '
$ws.Range("D27").Value = '''135.16'
$ws.Range("E27").Value = '  +1.86%  '
$ws.Range("D28").Value = '''17.05'
$ws.Range("E28").Value = '  +3.08%  '
$ws.Range("D29").Value = '''1.990'
$ws.Range("E29").Value = '  +5.21%  '
$ws.Range("D30").Value = '''1.391'
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("D31").Value = '''4.221'
$ws.Range("E31").Value = '  +1.19%  '
$ws.Range("D32").Value = '''0.08781'
$ws.Range("E32").Value = '  +4.90%  '
$ws.Range("D33").Value = '''4.000'
$ws.Range("E33").Value = '  +3.73%  '
$ws.Range("D34").Value = '''0.05058'
$ws.Range("E34").Value = '  +2.14%  '
$ws.Range("D35").Value = '''1.142'
$ws.Range("E35").Value = '  +4.53%  '
$ws.Range("D36").Value = '''0.7122'
$ws.Range("E36").Value = '  +4.74%  '
$ws.Range("E37").Value = '  -0.69%  '
$ws.Range("D38").Value = '''2.809'
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("D39").Value = '''2.268'
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("D40").Value = '''0.9689'
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").Value = '''0.01692'
$ws.Range("E41").Value = '  +6.15%  '
$ws.Range("D42").Value = '''6.127'
$ws.Range("E42").Value = '  +0.72%  '
$ws.Range("D43").Value = '''0.4282'
$ws.Range("E43").Value = '  +4.30%  '
$ws.Range("D44").Value = '''104.62'
$ws.Range("E44").Value = '  +3.84%  '
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("D46").Value = '''7.575'
$ws.Range("E46").Value = '  +4.85%  '
$ws.Range("D47").Value = '''0.1275'
$ws.Range("E47").Value = '  +4.05%  '
$ws.Range("D48").Value = '''0.05737'
$ws.Range("E48").Value = '  +4.06%  '
$ws.Range("D49").Value = '''33.05'
$ws.Range("E49").Value = '  +4.20%  '
$ws.Range("D50").Value = '''8.436'
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("D51").Value = '''0.3791'
$ws.Range("E51").Value = '  +4.33%  '
